$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.891.60'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.094.05'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.35'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.17'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.088.47'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.50'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.156'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.58%  '

$ws.Range("E13").Value = '  +3.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.84'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.594.72'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.023.68'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.47%  '

$ws.Range("E17").Value = '  +1.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.090.71'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.66'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.82'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.45'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.700'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.10'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.92'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.36'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("E27").Value = '  -0.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.06'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.93%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.35'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.90'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.17%  '

$ws.Range("E32").Value = '  +1.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '57.31'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.36'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.41'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +5.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '497.90'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.02'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.246.40'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0404'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0800'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.70'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.13'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.254'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '124.31'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.26'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.04'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0533'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.11%  '

$ws.Range("E50").Value = '  +1.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.34'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.06%  '
